$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Report generation for handback: populate the "Latest Target File",
# "Latest Handback File" and "Latest Handback DateTime" columns on the
# per-locale sheets, refresh the overall Status text, and widen the columns
# that now hold longer content.
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# Widen the (now longer) status columns on the Overview sheet.
$wsOverview.Cells.Item(1, 5).EntireColumn.ColumnWidth = 29.166666666666668
$wsOverview.Cells.Item(1, 6).EntireColumn.ColumnWidth = 29.166666666666668

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText

$wsZh.Hyperlinks.Add(
    $wsZh.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c918f4590f4bd4cfd6c632b6f7b238f8ffd2121/e2e/1a99bda7-55d2-4d90-9662-5ba4908eab92.md",
    [Type]::Missing,
    [Type]::Missing,
    "1a99bda7-55d2-4d90-9662-5ba4908eab92.md"
)
$wsZh.Range("J2").Value = "1a99bda7-55d2-4d90-9662-5ba4908eab92.66861ffe82cf035077894158410bd0cbf7567fa9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 11:04:58"

$wsZh.Cells.Item(1, 3).EntireColumn.ColumnWidth = 29.166666666666668
$wsZh.Cells.Item(1, 9).EntireColumn.ColumnWidth = 39.166666666666664
$wsZh.Cells.Item(1, 10).EntireColumn.ColumnWidth = 39.166666666666664

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText

$wsDe.Hyperlinks.Add(
    $wsDe.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6c918f4590f4bd4cfd6c632b6f7b238f8ffd2121/e2e/1a99bda7-55d2-4d90-9662-5ba4908eab92.md",
    [Type]::Missing,
    [Type]::Missing,
    "1a99bda7-55d2-4d90-9662-5ba4908eab92.md"
)
$wsDe.Range("J2").Value = "1a99bda7-55d2-4d90-9662-5ba4908eab92.66861ffe82cf035077894158410bd0cbf7567fa9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 11:05:16"

$wsDe.Cells.Item(1, 3).EntireColumn.ColumnWidth = 29.166666666666668
$wsDe.Cells.Item(1, 9).EntireColumn.ColumnWidth = 39.166666666666664
$wsDe.Cells.Item(1, 10).EntireColumn.ColumnWidth = 39.166666666666664
